$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New coordinate rows (rows 2-10), appended below the existing two rows.
$data = @(
    @(40.73061, -73.935242),
    @(35.7721, -78.63861),
    @(38.123, -78.543),
    @(35.60096, -81.6467086666667),
    @(34.297155, -83.9505876666667),
    @(32.99335, -86.2544666666667),
    @(31.689545, -88.5583456666667),
    @(30.38574, -90.8622246666667),
    @(29.081935, -93.1661036666667)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

# Build the combined number-format + alignment style on a scratch cell
# outside the target range so the data range only ever picks up the
# single final style (avoids leaving an extra unused intermediate style
# behind in cellXfs).
$tmpl = $ws.Range("D1")
$tmpl.NumberFormat = "@"
$tmpl.HorizontalAlignment = -4131
$tmpl.Copy()
$ws.Range("A2:B10").PasteSpecial(-4122)
$tmpl.Clear()
$excel.CutCopyMode = $false

# Selection moves to B8 after the edits.
$ws.Range("B8").Select()

$ws.PageSetup.Orientation = 1
